# Generate Report for Handoff
# Update localization-status report: the 9e45521e-a2d2-475d-85ea-8856077458eb.md
# file moved from "Handed back: in sync with en-US" to "Ready for handoff" for
# both zh-cn and de-de locales, with refreshed handoff timestamps and an
# error detail explaining the stale handback version.

$wb = $excel.ActiveWorkbook

$statusReadyForHandoff = "Ready for handoff"
$overviewDate = "2016-08-13 17:05:43"
$zhcnHandoffDate = "2016-08-13 17:05:35"
$dedeHandoffDate = "2016-08-13 17:05:43"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/3a2079949632e16255b52f0ec3382f0b0a2bbd76/e2e/9e45521e-a2d2-475d-85ea-8856077458eb.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/b2b98dc7c80b0446bb78ad46da25eb30d60ccc03/e2e/9e45521e-a2d2-475d-85ea-8856077458eb.md."

# --- Overview sheet: row for 9e45521e-...md (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusReadyForHandoff
$wsOverview.Range("F3").Value = $statusReadyForHandoff
$wsOverview.Range("G3").Value = $overviewDate

# --- zh-cn sheet: row for 9e45521e-...md (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusReadyForHandoff
$wsZhCn.Range("H3").Value = $zhcnHandoffDate
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.14285714285714

# --- de-de sheet: row for 9e45521e-...md (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusReadyForHandoff
$wsDeDe.Range("H3").Value = $dedeHandoffDate
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.14285714285714
